$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3788.238
$ws.Range("J17").Value = 3877.65
$ws.Range("L17").Value = 11632.95
$ws.Range("N17").Value = -11968.95
$ws.Range("H70").Value = 5252.4165
$ws.Range("I70").Value = 849.5
$ws.Range("K70").Value = 2548.5
$ws.Range("M70").Value = -2278.5
$ws.Range("H73").Value = 5252.4165
$ws.Range("I73").Value = 849.5
$ws.Range("K73").Value = 2548.5
$ws.Range("M73").Value = -1612.5
$ws.Range("H106").Value = 2573.5833
$ws.Range("I106").Value = 2398.5557
$ws.Range("J106").Value = 3098.6667
$ws.Range("K106").Value = 2398.5557
$ws.Range("L106").Value = 3098.6667
$ws.Range("M106").Value = -1767.5557
$ws.Range("N106").Value = -4360.6667
$ws.Range("H107").Value = 2475
$ws.Range("I107").Value = 2475
$ws.Range("K107").Value = 2475
$ws.Range("M107").Value = -555
$ws.Range("H113").Value = 2145.0513
$ws.Range("I113").Value = 1958.0312
$ws.Range("K113").Value = 1958.0312
$ws.Range("M113").Value = 1295.9688
$ws.Range("H130").Value = 59985
$ws.Range("J130").Value = 59985
$ws.Range("L130").Value = 59985
$ws.Range("N130").Value = -70025
$ws.Range("H140").Value = 124999
$ws.Range("J140").Value = 124999
$ws.Range("L140").Value = 124999
$ws.Range("N140").Value = -135359

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 288
$ws.Range("J76").Value = 288
$ws.Range("L76").Value = 288
$ws.Range("N76").Value = -964
$ws.Range("H79").Value = 288
$ws.Range("J79").Value = 288
$ws.Range("L79").Value = 288
$ws.Range("N79").Value = -2628
$ws.Range("H110").Value = 1789.2
$ws.Range("I110").Value = 1789.2
$ws.Range("K110").Value = 1789.2
$ws.Range("M110").Value = 255.8
$ws.Range("H125").Value = 31000
$ws.Range("J125").Value = 31000
$ws.Range("L125").Value = 31000
$ws.Range("N125").Value = -40840

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 5001.5
$ws.Range("I105").Value = 5001.5
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 5001.5
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -3254.5
$ws.Range("N105").ClearContents()
$ws.Range("H107").Value = 1244.1818
$ws.Range("I107").Value = 810.1667
$ws.Range("J107").Value = 1765
$ws.Range("K107").Value = 810.1667
$ws.Range("L107").Value = 1765
$ws.Range("M107").Value = 1109.8333
$ws.Range("N107").Value = -5605
$ws.Range("H134").Value = 4716.4736
$ws.Range("I134").Value = 2976.5757
$ws.Range("J134").Value = 16199.8
$ws.Range("K134").Value = 8929.7271
$ws.Range("L134").Value = 48599.39999999999
$ws.Range("M134").Value = -6394.7271
$ws.Range("N134").Value = -53669.39999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1450
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H31").Value = 8069.7295
$ws.Range("I31").Value = 4014.2632
$ws.Range("J31").Value = 12350.5
$ws.Range("K31").Value = 4014.2632
$ws.Range("L31").Value = 12350.5
$ws.Range("M31").Value = -3719.2632
$ws.Range("N31").Value = -12940.5
$ws.Range("H34").Value = 8069.7295
$ws.Range("I34").Value = 4014.2632
$ws.Range("J34").Value = 12350.5
$ws.Range("K34").Value = 4014.2632
$ws.Range("L34").Value = 12350.5
$ws.Range("M34").Value = -3812.2632
$ws.Range("N34").Value = -12754.5
$ws.Range("H58").Value = 4405.4165
$ws.Range("I58").Value = 5224.75
$ws.Range("J58").Value = 3995.75
$ws.Range("K58").Value = 5224.75
$ws.Range("L58").Value = 3995.75
$ws.Range("M58").Value = -5021.75
$ws.Range("N58").Value = -4401.75
$ws.Range("H74").Value = 28227.4
$ws.Range("I74").Value = 10284.5
$ws.Range("K74").Value = 10284.5
$ws.Range("M74").Value = -9410.5
$ws.Range("H77").Value = 28227.4
$ws.Range("I77").Value = 10284.5
$ws.Range("K77").Value = 30853.5
$ws.Range("M77").Value = -26485.5
$ws.Range("H99").Value = 3914.9092
$ws.Range("I99").Value = 3895.875
$ws.Range("J99").Value = 3965.6667
$ws.Range("K99").Value = 3895.875
$ws.Range("L99").Value = 3965.6667
$ws.Range("M99").Value = -2397.875
$ws.Range("N99").Value = -6961.6667
$ws.Range("H104").Value = 34285
$ws.Range("J104").Value = 34285
$ws.Range("L104").Value = 34285
$ws.Range("N104").Value = -39527
$ws.Range("H113").Value = 1450
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H126").Value = 3914.9092
$ws.Range("I126").Value = 3895.875
$ws.Range("J126").Value = 3965.6667
$ws.Range("K126").Value = 11687.625
$ws.Range("L126").Value = 11897.0001
$ws.Range("M126").Value = -9217.625
$ws.Range("N126").Value = -16837.0001
$ws.Range("H132").Value = 7258.8823
$ws.Range("I132").Value = 5108.8335
$ws.Range("K132").Value = 15326.5005
$ws.Range("M132").Value = -12796.5005
$ws.Range("H134").Value = 4154.579
$ws.Range("I134").Value = 3919.9167
$ws.Range("J134").Value = 4556.857
$ws.Range("K134").Value = 11759.7501
$ws.Range("L134").Value = 13670.571
$ws.Range("M134").Value = -9224.750100000001
$ws.Range("N134").Value = -18740.571
$ws.Range("H136").Value = 4405.4165
$ws.Range("I136").Value = 5224.75
$ws.Range("J136").Value = 3995.75
$ws.Range("K136").Value = 15674.25
$ws.Range("L136").Value = 11987.25
$ws.Range("M136").Value = -13124.25
$ws.Range("N136").Value = -17087.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 562.9167
$ws.Range("I8").Value = 562.9167
$ws.Range("K8").Value = 1688.7501
$ws.Range("M8").Value = -1549.7501
$ws.Range("H47").Value = 314.33334
$ws.Range("I47").Value = 314.33334
$ws.Range("K47").Value = 943.0000200000001
$ws.Range("M47").Value = -512.0000200000001
$ws.Range("H107").Value = 420.83334
$ws.Range("J107").Value = 445.1111
$ws.Range("L107").Value = 1335.3333
$ws.Range("N107").Value = -5175.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 48671
$ws.Range("J104").Value = 48671
$ws.Range("L104").Value = 48671
$ws.Range("N104").Value = -55659
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("H113").Value = 3861.75
$ws.Range("I113").Value = 2750
$ws.Range("K113").Value = 2750
$ws.Range("M113").Value = -580
$ws.Range("H132").Value = 6537.8335
$ws.Range("I132").Value = 4023.5
$ws.Range("J132").Value = 10057.9
$ws.Range("K132").Value = 12070.5
$ws.Range("L132").Value = 30173.7
$ws.Range("M132").Value = -9540.5
$ws.Range("N132").Value = -35233.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H51").Value = 30000
$ws.Range("J51").Value = 30000
$ws.Range("L51").Value = 30000
$ws.Range("N51").Value = -30956
$ws.Range("H61").Value = 6761.375
$ws.Range("J61").Value = 9484.6
$ws.Range("L61").Value = 9484.6
$ws.Range("N61").Value = -9888.6
$ws.Range("H113").Value = 6761.375
$ws.Range("J113").Value = 9484.6
$ws.Range("L113").Value = 9484.6
$ws.Range("N113").Value = -13824.6
$ws.Range("H127").Value = 88832.664
$ws.Range("J127").Value = 88832.664
$ws.Range("L127").Value = 88832.664
$ws.Range("N127").Value = -98752.664
$ws.Range("H136").Value = 7550.0454
$ws.Range("I136").Value = 5106.1665
$ws.Range("J136").Value = 9241.962
$ws.Range("K136").Value = 15318.4995
$ws.Range("L136").Value = 27725.886
$ws.Range("M136").Value = -12768.4995
$ws.Range("N136").Value = -32825.886

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 120000
$ws.Range("J46").Value = 120000
$ws.Range("L46").Value = 120000
$ws.Range("N46").Value = -120462
$ws.Range("H62").Value = 10996
$ws.Range("I62").Value = 10002
$ws.Range("J62").Value = 11990
$ws.Range("K62").Value = 10002
$ws.Range("L62").Value = 11990
$ws.Range("M62").Value = -9378
$ws.Range("N62").Value = -13238
$ws.Range("H65").Value = 10996
$ws.Range("I65").Value = 10002
$ws.Range("J65").Value = 11990
$ws.Range("K65").Value = 50010
$ws.Range("L65").Value = 59950
$ws.Range("M65").Value = -46890
$ws.Range("N65").Value = -66190
$ws.Range("H107").Value = 1798.4482
$ws.Range("I107").Value = 2222.125
$ws.Range("J107").Value = 1277
$ws.Range("K107").Value = 6666.375
$ws.Range("L107").Value = 3831
$ws.Range("M107").Value = -4746.375
$ws.Range("N107").Value = -7671
$ws.Range("H134").Value = 120000
$ws.Range("J134").Value = 120000
$ws.Range("L134").Value = 360000
$ws.Range("N134").Value = -365070
$ws.Range("H136").Value = 3456.8064
$ws.Range("I136").Value = 2972.0334
$ws.Range("J136").Value = 18000
$ws.Range("K136").Value = 8916.100199999999
$ws.Range("L136").Value = 54000
$ws.Range("M136").Value = -6366.100199999999
$ws.Range("N136").Value = -59100
